$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in a new row of data (Day 23) in the previously-blank row 26,
# just above the "Total hours Spent" block (which stays put at row 36).
$ws.Cells.Item(26, 1).Value = 23
$ws.Cells.Item(26, 2).Value = "27/3/2024"
$ws.Cells.Item(26, 3).Value = 3.75
$ws.Cells.Item(26, 4).Value = "Made Product Update functionality + a lot of refactoring"

# Match formatting (center aligned) used by the rest of the data rows.
$ws.Range("A26:D26").HorizontalAlignment = -4108
$ws.Range("A26:D26").VerticalAlignment = -4108

# Update the active selection to reflect where the user left off editing.
$ws.Range("D29").Select()

$wb.Save()
